# The Tracability Link Matrix had three columns removed: "GameDriver",
# "GameHistory" and "LeaderBoard" (originally columns D, F and G). Deleting
# these entire columns shifts the remaining columns (Game, Invitation, User,
# DBconnection) left, which is exactly what is reflected in the target
# worksheet (now only using columns A:G instead of A:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from right to left so earlier deletions don't shift the
# not-yet-deleted column letters out from under us.
$ws.Range("G1").EntireColumn.Delete()   # LeaderBoard
$ws.Range("F1").EntireColumn.Delete()   # GameHistory
$ws.Range("D1").EntireColumn.Delete()   # GameDriver

# Leave the same selection state captured in the saved workbook: the whole
# of (the new) column D -- which now holds "Game" -- selected.
$ws.Range("D1").EntireColumn.Select() | Out-Null
